$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.256905436515808
$ws.Range("B1").Value = 2.686847925186157
$ws.Range("C1").Value = 5.043686389923096
$ws.Range("D1").Value = 2.061823129653931
$ws.Range("E1").Value = 1.038701772689819
